$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text, matching the
# original inline-string cell type, so numeric-looking values like
# "21.50" or "0.780" keep their exact formatting instead of being
# coerced into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.043.59'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("D3").Value = '1.680.68'
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '215.94'
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("E6").Value = '  -3.26%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '0.255'
$ws.Range("E8").Value = '  +1.42%  '

$ws.Range("D9").Value = '21.50'
$ws.Range("E9").Value = '  +6.67%  '

$ws.Range("D10").Value = '0.0624'
$ws.Range("E10").Value = '  +0.92%  '

$ws.Range("E11").Value = '  -0.77%  '

$ws.Range("D12").Value = '1.917.89'
$ws.Range("E12").Value = '  +0.90%  '

$ws.Range("D13").Value = '1.677.63'
$ws.Range("E13").Value = '  -0.09%  '

$ws.Range("E14").Value = '  +0.70%  '

$ws.Range("E15").Value = '  +1.75%  '

$ws.Range("D16").Value = '66.52'
$ws.Range("E16").Value = '  +0.83%  '

$ws.Range("D17").Value = '27.041.74'
$ws.Range("E17").Value = '  +0.65%  '

$ws.Range("E18").Value = '  +4.96%  '

$ws.Range("D19").Value = '235.76'
$ws.Range("E19").Value = '  +1.80%  '

$ws.Range("D20").Value = '0.0₃0740'
$ws.Range("E20").Value = '  +1.47%  '

$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("E22").Value = '  +0.57%  '

$ws.Range("D23").Value = '9.28'
$ws.Range("E23").Value = '  +1.32%  '

$ws.Range("E24").Value = '  -4.24%  '

$ws.Range("D25").Value = '146.65'
$ws.Range("E25").Value = '  +0.70%  '

$ws.Range("D26").Value = '7.24'
$ws.Range("E26").Value = '  +1.81%  '

$ws.Range("D27").Value = '16.49'
$ws.Range("E27").Value = '  +4.03%  '

$ws.Range("E28").Value = '  -2.57%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  +0.63%  '

$ws.Range("D31").Value = '1.18'
$ws.Range("E31").Value = '  +0.33%  '

$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  -0.24%  '

$ws.Range("D33").Value = '1.529.79'
$ws.Range("E33").Value = '  +4.58%  '

$ws.Range("E34").Value = '  +0.97%  '

$ws.Range("E35").Value = '  +5.25%  '

$ws.Range("E36").Value = '  -0.79%  '

$ws.Range("D37").Value = '0.591'
$ws.Range("E37").Value = '  +3.21%  '

$ws.Range("E38").Value = '  +2.81%  '

$ws.Range("E39").Value = '  +3.34%  '

$ws.Range("E40").Value = '  +5.28%  '

$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").Value = '5.69'
$ws.Range("E42").Value = '  -3.01%  '

$ws.Range("D43").Value = '67.88'
$ws.Range("E43").Value = '  +3.18%  '

$ws.Range("D44").Value = '2.25'
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("D45").Value = '1.822.90'
$ws.Range("E45").Value = '  +0.47%  '

$ws.Range("D46").Value = '0.780'
$ws.Range("E46").Value = '  +0.49%  '

$ws.Range("D47").Value = '90.22'
$ws.Range("E47").Value = '  +0.00%  '

$ws.Range("E48").Value = '  -0.08%  '

$ws.Range("E49").Value = '  +2.69%  '

$ws.Range("D50").Value = '7.99'
$ws.Range("E50").Value = '  +5.80%  '

$ws.Range("D51").Value = '0.0504'
$ws.Range("E51").Value = '  -0.42%  '
